$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 501, shifting existing rows 501-605 down to 502-606.
$ws.Rows("501").Insert()

$ws.Range("A501").Value = 6
$ws.Range("B501").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C501").Value = "Metropolitana"
$ws.Range("D501").Value = 44943
$ws.Range("E501").Value = 13
$ws.Range("F501").Value = 100112030
$ws.Range("G501").Value = "Poroto granado"
$ws.Range("H501").Value = "Sin especificar"
$ws.Range("I501").Value = "Primera"
$ws.Range("J501").Value = 1070
$ws.Range("K501").Value = 43000
$ws.Range("L501").Value = 45000
$ws.Range("M501").Value = 44215
$ws.Range("N501").Value = "`$/saco 25 kilos"
$ws.Range("O501").Value = "Región Metropolitana"
$ws.Range("P501").Value = 1769
$ws.Range("Q501").Value = 25
$ws.Range("R501").Value = "Hortaliza"
